$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the thick-bottom-border band into columns Q:R (and restyle P3) ---
$ws.Range("O3").Copy()
$ws.Range("P3:R3").PasteSpecial(-4122)

# --- Build the "year header" style (bold 10pt Times New Roman, medium bottom border) ---
# on a scratch cell, then stamp it onto P4:R4 and fill in the new years.
$ws.Range("A8").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("Z100").Font.Bold = $true
$ws.Range("Z100").Font.Size = 10
$ws.Range("Z100").Borders.Item(9).LineStyle = 1
$ws.Range("Z100").Borders.Item(9).Weight = -4138
$ws.Range("Z100").Borders.Item(9).ColorIndex = 1
$ws.Range("Z100").Copy()
$ws.Range("P4:R4").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

$ws.Range("P4").Value = 2019
$ws.Range("Q4").Value = 2020
$ws.Range("R4").Value = 2021

# --- Build the "data value" style (10pt Times New Roman, top+bottom medium border)
# on a scratch cell, then stamp it onto P5:R5 and fill in the new figures.
$ws.Range("E5").Copy()
$ws.Range("Z101").PasteSpecial(-4122)
$ws.Range("Z101").Font.Size = 10
$ws.Range("Z101").Copy()
$ws.Range("P5:R5").PasteSpecial(-4122)
$ws.Range("Z101").Clear()

$ws.Range("P5").Value = 12.9
$ws.Range("Q5").Value = 15.2
$ws.Range("R5").Value = 10.4

# --- Update the active-cell selection to match the authored state ---
$ws.Range("S3").Select()
